# Regenerate save_data: use K (column G) instead of Strike# values,
# recomputed from refreshed std/mean calc of s_vals.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column G ("K")
$newK = @{
    2  = 2
    3  = 0
    4  = 1
    5  = 0
    6  = 1
    7  = 1
    8  = 1
    9  = 1
    10 = 1
    11 = 0
    13 = 1
    14 = 0
    15 = 1
    16 = 2
    17 = 0
    18 = 1
    19 = 0
    20 = 0
    21 = 1
    22 = 0
    23 = 1
    24 = 0
    25 = 1
    26 = 0
    27 = 2
    28 = 1
    29 = 1
    30 = 0
    31 = 0
    32 = 1
    33 = 2
    34 = 2
    35 = 0
    36 = 0
    37 = 0
    38 = 1
    39 = 1
    40 = 1
    41 = 1
    42 = 0
    43 = 2
    44 = 1
    45 = 1
    46 = 1
    47 = 1
    48 = 1
    49 = 1
    50 = 2
    51 = 1
    52 = 1
    53 = 2
    54 = 3
    55 = 0
    56 = 1
    57 = 1
    58 = 0
    59 = 1
    60 = 0
    61 = 0
    62 = 3
    63 = 2
    64 = 2
    65 = 2
    66 = 0
    67 = 1
    68 = 3
    70 = 2
    71 = 1
    72 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
